# Import / Export XL Fixes #25
#
# Adds a new "ID" column as the first column of each data sheet
# (Aclass, Bclass, Dclass), shifting the previously-existing columns
# one position to the right, and populates it with sequential record
# identifiers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Aclass": A1:L3 (Name..Duration1) -> A1:M3 (ID, Name..Duration1)
# ---------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Aclass")

# Wipe any inherited number formatting so the rebuilt grid starts clean.
$wsA.Range("A1:M3").ClearFormats()

$aHeaders = @("ID","Name","Date","Booleanfield","Aenum","Aenum_2","Benum","CName","CFloatfield","Floatfield","Intfield","Anotherbooleanfield","Duration1")
for ($c = 1; $c -le $aHeaders.Length; $c++) {
    $wsA.Cells.Item(1, $c).Value = $aHeaders[$c - 1]
}

# Clear the stale, pre-shift contents of the whole data block first so
# no leftover values from the old column layout survive in cells that
# land on blank ("") columns below (CName / Aenum_2 / Benum).
$wsA.Range("A2:M3").ClearContents()

$wsA.Cells.Item(2, 1).Value = 1
$wsA.Cells.Item(2, 2).Value = "A1"
$wsA.Cells.Item(2, 3).Value = 44247.84171296297
$wsA.Cells.Item(2, 3).NumberFormat = "m/d/yy h:mm"
$wsA.Cells.Item(2, 4).Value = $true
$wsA.Cells.Item(2, 5).Value = "ENUM_VAL1_NOT_THE_SAME"
# columns 6 (Aenum_2), 7 (Benum), 8 (CName) stay blank
$wsA.Cells.Item(2, 9).Value = 0
$wsA.Cells.Item(2, 10).Value = 10.2
$wsA.Cells.Item(2, 11).Value = 4
$wsA.Cells.Item(2, 12).Value = $true
$wsA.Cells.Item(2, 13).Value = "1h3m0.001s"

$wsA.Cells.Item(3, 1).Value = 3
$wsA.Cells.Item(3, 2).Value = "A2"
$wsA.Cells.Item(3, 3).Value = -693593
$wsA.Cells.Item(3, 3).NumberFormat = "m/d/yy h:mm"
$wsA.Cells.Item(3, 4).Value = $true
# columns 5 (Aenum), 6 (Aenum_2), 7 (Benum), 8 (CName) stay blank
$wsA.Cells.Item(3, 9).Value = 0
$wsA.Cells.Item(3, 10).Value = 10.77
$wsA.Cells.Item(3, 11).Value = 0
$wsA.Cells.Item(3, 12).Value = $true
$wsA.Cells.Item(3, 13).Value = "0s"

# ---------------------------------------------------------------------
# Sheet "Bclass": A1:C3 (Name, Floatfield, Intfield) -> A1:D3 (ID, Name, Floatfield, Intfield)
# ---------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("Bclass")

$wsB.Range("A1:D3").ClearFormats()

$bHeaders = @("ID","Name","Floatfield","Intfield")
for ($c = 1; $c -le $bHeaders.Length; $c++) {
    $wsB.Cells.Item(1, $c).Value = $bHeaders[$c - 1]
}

$wsB.Cells.Item(2, 1).Value = 1
$wsB.Cells.Item(2, 2).Value = "B1"
$wsB.Cells.Item(2, 3).Value = 0
$wsB.Cells.Item(2, 4).Value = 0

$wsB.Cells.Item(3, 1).Value = 2
$wsB.Cells.Item(3, 2).Value = "B2"
$wsB.Cells.Item(3, 3).Value = 0
$wsB.Cells.Item(3, 4).Value = 0

# ---------------------------------------------------------------------
# Sheet "Dclass": A1 (Name) -> A1:B1 (ID, Name)
# ---------------------------------------------------------------------
$wsD = $wb.Worksheets.Item("Dclass")

$wsD.Range("A1:B1").ClearFormats()

$wsD.Cells.Item(1, 1).Value = "ID"
$wsD.Cells.Item(1, 2).Value = "Name"
